$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '64.283.84'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = "'" + '3.397.56'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'" + '568.31'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").Value = "'" + '156.20'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("D7").Value = "'" + '0.632'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +9.52%  '
$ws.Range("D8").Value = "'" + '1.00'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = "'" + '3.407.19'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("D10").Value = "'" + '7.13'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D11").Value = "'" + '0.122'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("D12").Value = "'" + '0.439'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.83%  '
$ws.Range("D13").Value = "'" + '3.984.50'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = "'" + '0.134'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = "'" + '0.0000187'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.53%  '
$ws.Range("D16").Value = "'" + '27.43'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("D17").Value = "'" + '64.282.95'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = "'" + '3.445.69'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("D19").Value = "'" + '6.28'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").Value = "'" + '13.75'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.67%  '
$ws.Range("D21").Value = "'" + '376.78'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.81%  '
$ws.Range("D22").Value = "'" + '7.97'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("D23").Value = "'" + '0.543'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").Value = "'" + '0.999'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = "'" + '71.73'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.35%  '
$ws.Range("D26").Value = "'" + '0.0000118'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.70%  '
$ws.Range("D27").Value = "'" + '10.26'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +7.84%  '
$ws.Range("D28").Value = "'" + '0.178'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("D29").Value = "'" + '0.999'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").Value = "'" + '1.46'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.19%  '
$ws.Range("D31").Value = "'" + '6.13'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").Value = "'" + '2.00'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.45%  '
$ws.Range("D33").Value = "'" + '22.97'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("D34").Value = "'" + '7.12'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.87%  '
$ws.Range("D35").Value = "'" + '1.60'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +7.72%  '
$ws.Range("D36").Value = "'" + '159.47'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("D37").Value = "'" + '1.89'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("D38").Value = "'" + '6.93'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +6.16%  '
$ws.Range("D39").Value = "'" + '0.0757'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.41%  '
$ws.Range("D40").Value = "'" + '2.875.47'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.62%  '
$ws.Range("D41").Value = "'" + '4.61'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.39%  '
$ws.Range("D42").Value = "'" + '26.17'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("D43").Value = "'" + '42.89'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("D44").Value = "'" + '0.0313'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("D45").Value = "'" + '25.88'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.33%  '
$ws.Range("D46").Value = "'" + '0.766'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = "'" + '321.01'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +7.42%  '
$ws.Range("D48").Value = "'" + '1.07'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = "'" + '0.110'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.31%  '
$ws.Range("D50").Value = "'" + '2.19'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.59%  '
$ws.Range("D51").Value = "'" + '6.54'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.78%  '
